# Weekly update: two new "Cebollín" price records were reported for
# Terminal Hortofrutícola Agro Chillán, dated 2022-11-11 (serial 44876).
# They are inserted as the new first two data rows (worksheet rows 32-33),
# pushing the existing rows 32-38 down to rows 34-40 unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 32, shifting rows 32-38
# (and everything below) down by two rows.
$ws.Rows("32:33").Insert()

# --- New row 32: "Primera" quality, $/paquete 6 unidades ---
$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C32").Value = "Ñuble"
$ws.Range("D32").Value = 44876
$ws.Range("E32").Value = 16
$ws.Range("F32").Value = 100112037
$ws.Range("G32").Value = "Cebollín"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 500
$ws.Range("K32").Value = 600
$ws.Range("L32").Value = 700
$ws.Range("M32").Value = 650
$ws.Range("N32").Value = "`$/paquete 6 unidades"
$ws.Range("O32").Value = "Provincia de Diguillín"
$ws.Range("P32").Value = 108
$ws.Range("Q32").Value = 6
$ws.Range("R32").Value = "Hortaliza"

# --- New row 33: "Segunda" quality, $/paquete 6 unidades ---
$ws.Range("A33").Value = 7
$ws.Range("B33").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C33").Value = "Ñuble"
$ws.Range("D33").Value = 44876
$ws.Range("E33").Value = 16
$ws.Range("F33").Value = 100112037
$ws.Range("G33").Value = "Cebollín"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Segunda"
$ws.Range("J33").Value = 400
$ws.Range("K33").Value = 500
$ws.Range("L33").Value = 500
$ws.Range("M33").Value = 500
$ws.Range("N33").Value = "`$/paquete 6 unidades"
$ws.Range("O33").Value = "Provincia de Diguillín"
$ws.Range("P33").Value = 83
$ws.Range("Q33").Value = 6
$ws.Range("R33").Value = "Hortaliza"
